$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet SRS_dict -> ART_dict
$ws.Name = "ART_dict"

# --- Row 3 (INSTRUCTIONS): update German text (ein Literat -> eine Literatin, Dramatiker:in added) ---
$ws.Range("B3").Value = "Sie werden nun mehrere Paare von Namen sehen und müssen jeweils entscheiden, welcher der Persoen eine Literatin (Dichter:in, Schriftsteller:in, Dramatiker:in) ist."

# --- Row 4 (PROMPT): update German + English text ---
$ws.Range("B4").Value = "Welcher der Personen ist eine Literat:in?<br/> Klicken Sie auf den Namen, sie haben {{time_out}} Sekunden Zeit zu antworten."
$ws.Range("C4").Value = "Who of the persons is a literary writer? <br/> Click on the name, you have {{time_out}} seconds."

# --- Insert two new rows (PROMPT_SINGLE, PROMPT_SINGLE_PAGE) before the old FEEDBACK row ---
$ws.Rows("5:6").Insert() | Out-Null
$ws.Range("A5").Value = "PROMPT_SINGLE"
$ws.Range("B5").Value = "Ist <b>{{name}}</b> eine Literat:in?<br/> Klicken Sie Ja oder Nein, sie haben {{time_out}} Sekunden Zeit zu antworten."
$ws.Range("C5").Value = "Is <b>{{name}}</b> a literary writer? <br/> Click Yes or No, you have {{time_out}} seconds."
$ws.Range("A6").Value = "PROMPT_SINGLE_PAGE"
$ws.Range("B6").Value = "Bitte wählen Sie alle Literat:innen (Dichter:innen, Romanautor:innen, Dramatiker:innen)  aus der untenstehenden Liste aus.  Sie haben {{time_out}} Sekunden Zeit."
$ws.Range("C6").Value = "Please select all literary writers (poets, novelists, playwrights). You have {{time_out}} seconds."

# --- Row 7 is now FEEDBACK (text unchanged); insert a new row after it for FEEDBACK_SINGLE_PAGE ---
$ws.Rows("8:8").Insert() | Out-Null
$ws.Range("A8").Value = "FEEDBACK_SINGLE_PAGE"
$ws.Range("B8").Value = "Sie haben {{num_correct}} Literaten aus {{num_items}} Namen richtig erkannt ({{perc_correct}}%, Punkte: {{points}})."
$ws.Range("C8").Value = "You answered {{num_correct}} out of {{num_items}} questions correctly ({{perc_correct}}%,  Points: {{points}})."

# --- Append new YES/NO rows at the end ---
$ws.Range("A20").Value = "YES"
$ws.Range("B20").Value = "Ja"
$ws.Range("C20").Value = "Yes"
$ws.Range("A21").Value = "NO"
$ws.Range("B21").Value = "Nein"
$ws.Range("C21").Value = "No"

# Update selection to match the authored workbook
$ws.Range("C7").Select() | Out-Null
